# Natmi following Dr Hou advice
# Rebuild the LR-pairs (Pspn-Gfra1) results table with the updated
# cell-cell signalling rows: FAPs/M2 sending clusters paired with the
# Pspn-Gfra1 ligand-receptor pair targeting FAPs, ECs and sCs clusters.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Pspn"
$ws.Range("C2").Value = "Gfra1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 3.042572666666667
$ws.Range("H2").Value = 9.127718
$ws.Range("I2").Value = 0.9827375016055572
$ws.Range("J2").Value = 0.9827375016055572
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.509196
$ws.Range("N2").Value = 1.527588
$ws.Range("O2").Value = 0.02558190413389134
$ws.Range("P2").Value = 0.02558190413389134
$ws.Range("Q2").Value = 1.549265831576
$ws.Range("R2").Value = 13.943392484184
$ws.Range("S2").Value = 0.02514029655485325
$ws.Range("T2").Value = 0.02514029655485325

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Pspn"
$ws.Range("C3").Value = "Gfra1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 3.042572666666667
$ws.Range("H3").Value = 9.127718
$ws.Range("I3").Value = 0.9827375016055572
$ws.Range("J3").Value = 0.9827375016055572
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 17.676258
$ws.Range("N3").Value = 53.028774
$ws.Range("O3").Value = 0.8880516296316739
$ws.Range("P3").Value = 0.8880516296316739
$ws.Range("Q3").Value = 53.781299439748
$ws.Range("R3").Value = 484.031694957732
$ws.Range("S3").Value = 0.8727216398009748
$ws.Range("T3").Value = 0.8727216398009748

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Pspn"
$ws.Range("C4").Value = "Gfra1"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 3.042572666666667
$ws.Range("H4").Value = 9.127718
$ws.Range("I4").Value = 0.9827375016055572
$ws.Range("J4").Value = 0.9827375016055572
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.719084666666667
$ws.Range("N4").Value = 5.157254
$ws.Range("O4").Value = 0.0863664662344347
$ws.Range("P4").Value = 0.0863664662344347
$ws.Range("Q4").Value = 5.230440018485777
$ws.Range("R4").Value = 47.073960166372
$ws.Range("S4").Value = 0.08487556524972907
$ws.Range("T4").Value = 0.08487556524972907

# Row 5
$ws.Range("A5").Value = "M2"
$ws.Range("B5").Value = "Pspn"
$ws.Range("C5").Value = "Gfra1"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.053445
$ws.Range("H5").Value = 0.160335
$ws.Range("I5").Value = 0.01726249839444284
$ws.Range("J5").Value = 0.01726249839444284
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.509196
$ws.Range("N5").Value = 1.527588
$ws.Range("O5").Value = 0.02558190413389134
$ws.Range("P5").Value = 0.02558190413389134
$ws.Range("Q5").Value = 0.02721398022
$ws.Range("R5").Value = 0.24492582198
$ws.Range("S5").Value = 0.0004416075790380899
$ws.Range("T5").Value = 0.0004416075790380899

# Row 6
$ws.Range("A6").Value = "M2"
$ws.Range("B6").Value = "Pspn"
$ws.Range("C6").Value = "Gfra1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.053445
$ws.Range("H6").Value = 0.160335
$ws.Range("I6").Value = 0.01726249839444284
$ws.Range("J6").Value = 0.01726249839444284
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 17.676258
$ws.Range("N6").Value = 53.028774
$ws.Range("O6").Value = 0.8880516296316739
$ws.Range("P6").Value = 0.8880516296316739
$ws.Range("Q6").Value = 0.94470760881
$ws.Range("R6").Value = 8.50236847929
$ws.Range("S6").Value = 0.01532998983069912
$ws.Range("T6").Value = 0.01532998983069912

# Row 7
$ws.Range("A7").Value = "M2"
$ws.Range("B7").Value = "Pspn"
$ws.Range("C7").Value = "Gfra1"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.053445
$ws.Range("H7").Value = 0.160335
$ws.Range("I7").Value = 0.01726249839444284
$ws.Range("J7").Value = 0.01726249839444284
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.719084666666667
$ws.Range("N7").Value = 5.157254
$ws.Range("O7").Value = 0.0863664662344347
$ws.Range("P7").Value = 0.0863664662344347
$ws.Range("Q7").Value = 0.09187648001
$ws.Range("R7").Value = 0.82688832009
$ws.Range("S7").Value = 0.001490900984705631
$ws.Range("T7").Value = 0.001490900984705631
